$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated node table (Nodes, y, z) for rows 2-11 (10 data rows, was 6).
$rows = @(
    @(0, 0.5,   -0.5),
    @(1, -0.5,  -0.5),
    @(2, 0.25,  0.5),
    @(3, 0.25,  0.5),
    @(4, -0.25, 0.5),
    @(5, -0.25, 0.5),
    @(6, 0.25,  -0.5),
    @(7, -0.25, -0.5),
    @(8, 0,     0.5),
    @(9, 0,     -0.5)
)

$data = New-Object 'object[,]' $rows.Count,3
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 3; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("A2:C11").Value = $data

# Update the current selection to match the author's saved workbook state.
$ws.Range("G4:G6").Select()

$wb.Save()
